$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find and delete the rows for the cell types removed in this revision.
$namesToRemove = @("Endothelial cells", "Immune system cells", "Non myelinating Schwann cells")

foreach ($name in $namesToRemove) {
    $cell = $ws.Columns.Item(2).Find($name)
    if ($cell -ne $null) {
        $ws.Rows.Item($cell.Row).Delete()
    }
}

$ws.Range("B2").Select()
